# Fruta / hortaliza, semanal
#
# The weekly data refresh re-sequences the existing price observation rows
# (rows 2-39 on the single data sheet) into a new row order coming from the
# upstream export, without adding or removing any observations. Only the
# columns that vary per-observation move with each row (Fecha, Calidad,
# Volumen, Precio minimo/maximo/promedio ponderado, Origen, Precio $/Kg);
# the columns that are constant for every row in this subset (Mercado ID,
# Mercado, Region, Codreg, Tipo, Producto ID, Producto, Categoria ID,
# Categoria, Variedad, Unidad de comercializacion, Kg / unidad) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 39

# Column letter -> column index for the fields that are re-sequenced.
$cols = @{ "D" = 4; "L" = 12; "M" = 13; "N" = 14; "O" = 15; "P" = 16; "R" = 18; "S" = 19 }
# Columns that hold text (everything else in $cols is numeric).
$textCols = @("L", "R")

# New row number -> row number that currently holds the data that should end
# up there (i.e. destination row -> source row).
$rowMap = @{
    2=3;  3=10; 4=18;  5=32; 6=36;  7=11; 8=30;  9=21; 10=22; 11=4;
    12=35; 13=34; 14=5; 15=15; 16=37; 17=29; 18=23; 19=28; 20=14;
    21=7; 22=17; 23=38; 24=27; 25=31; 26=24; 27=12; 28=39; 29=16;
    30=26; 31=19; 32=9; 33=25; 34=6;  35=33; 36=13; 37=2;  38=8; 39=20
}

# 1) Snapshot every row's current values for the re-sequenced columns before
#    any writes happen (several rows are both a source and a destination).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols.Keys) {
        $colIndex = $cols[$col]
        $rowVals[$col] = $ws.Cells.Item($r, $colIndex).Formula
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each destination row using the snapshot captured from its mapped
#    source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols.Keys) {
        $colIndex = $cols[$col]
        $value = $srcVals[$col]
        if ($textCols -contains $col) {
            $ws.Cells.Item($destRow, $colIndex).Value = $value
        } else {
            $ws.Cells.Item($destRow, $colIndex).Value = [double]$value
        }
    }
}
